# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets. These two sheets contain identical data tables, and the
# same eight rows need their F-column figures bumped up slightly
# (re-scraped attendance counts).

$wb = $excel.ActiveWorkbook

$updates = @{
    7  = 1285
    8  = 1543
    10 = 405
    19 = 1744
    20 = 68
    28 = 278
    29 = 1100
    34 = 272
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
